$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Fgf7"
$ws.Range("C2").Value = "Fgfr3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 11.14069866666667
$ws.Range("H2").Value = 33.422096
$ws.Range("I2").Value = 0.9684602815214559
$ws.Range("J2").Value = 0.9684602815214559
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 2.619953333333334
$ws.Range("N2").Value = 7.85986
$ws.Range("O2").Value = 0.6100029538328192
$ws.Range("P2").Value = 0.6100029538328192
$ws.Range("Q2").Value = 29.18811060739556
$ws.Range("R2").Value = 262.69299546656
$ws.Range("S2").Value = 0.5907636323978518
$ws.Range("T2").Value = 0.5907636323978518

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Fgf7"
$ws.Range("C3").Value = "Fgfr3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 11.14069866666667
$ws.Range("H3").Value = 33.422096
$ws.Range("I3").Value = 0.9684602815214559
$ws.Range("J3").Value = 0.9684602815214559
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.5698483333333333
$ws.Range("N3").Value = 1.709545
$ws.Range("O3").Value = 0.1326776176306101
$ws.Range("P3").Value = 0.1326776176306101
$ws.Range("Q3").Value = 6.348508567368889
$ws.Range("R3").Value = 57.13657710632
$ws.Range("S3").Value = 0.1284930029221367
$ws.Range("T3").Value = 0.1284930029221367

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Fgf7"
$ws.Range("C4").Value = "Fgfr3"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 11.14069866666667
$ws.Range("H4").Value = 33.422096
$ws.Range("I4").Value = 0.9684602815214559
$ws.Range("J4").Value = 0.9684602815214559
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.105183
$ws.Range("N4").Value = 3.315549
$ws.Range("O4").Value = 0.2573194285365706
$ws.Range("P4").Value = 0.2573194285365706
$ws.Range("Q4").Value = 12.31251077452267
$ws.Range("R4").Value = 110.812596970704
$ws.Range("S4").Value = 0.2492036462014673
$ws.Range("T4").Value = 0.2492036462014673

# Row 5
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Fgf7"
$ws.Range("C5").Value = "Fgfr3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.3628176666666666
$ws.Range("H5").Value = 1.088453
$ws.Range("I5").Value = 0.03153971847854405
$ws.Range("J5").Value = 0.03153971847854405
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 2.619953333333334
$ws.Range("N5").Value = 7.85986
$ws.Range("O5").Value = 0.6100029538328192
$ws.Range("P5").Value = 0.6100029538328192
$ws.Range("Q5").Value = 0.9505653551755556
$ws.Range("R5").Value = 8.55508819658
$ws.Range("S5").Value = 0.01923932143496742
$ws.Range("T5").Value = 0.01923932143496742

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Fgf7"
$ws.Range("C6").Value = "Fgfr3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.3628176666666666
$ws.Range("H6").Value = 1.088453
$ws.Range("I6").Value = 0.03153971847854405
$ws.Range("J6").Value = 0.03153971847854405
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.5698483333333333
$ws.Range("N6").Value = 1.709545
$ws.Range("O6").Value = 0.1326776176306101
$ws.Range("P6").Value = 0.1326776176306101
$ws.Range("Q6").Value = 0.2067510426538889
$ws.Range("R6").Value = 1.860759383885
$ws.Range("S6").Value = 0.004184614708473354
$ws.Range("T6").Value = 0.004184614708473353

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Fgf7"
$ws.Range("C7").Value = "Fgfr3"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.3628176666666666
$ws.Range("H7").Value = 1.088453
$ws.Range("I7").Value = 0.03153971847854405
$ws.Range("J7").Value = 0.03153971847854405
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.105183
$ws.Range("N7").Value = 3.315549
$ws.Range("O7").Value = 0.2573194285365706
$ws.Range("P7").Value = 0.2573194285365706
$ws.Range("Q7").Value = 0.4009799172996666
$ws.Range("R7").Value = 3.608819255697
$ws.Range("S7").Value = 0.008115782335103273
$ws.Range("T7").Value = 0.008115782335103272
